$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 27 (Keyword "ZBook Firefly" / Model "ZBook Firefly"), shifting all
# subsequent rows up by one.
$ws.Rows.Item(27).Delete()

# Reset the view back to the top-left corner / default selection (A1) instead
# of leaving the scrolled-down view/selection from the editing session.
$ws.Range("A1").Select()
